# Append 5 new evaluation rows (24-28) to Sheet1, matching the rows that
# appear in the site's "evaluaciones_tk" export after the new submissions
# from jessica, marilyn, gypsi and luis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Evaluador, B Fiabilidad, C Eficiencia, D Mantenibilidad,
#          E Usabilidad, F Seguridad, G Fecha, H Promedio
$rows = @(
    @("jessica", 10, 10, 10, 10, 10, "2025-07-14 15:32:03", 10),
    @("marilyn",  5,  4,  3,  5,  5, "2025-07-14 15:34:02", 4.4),
    @("gypsi",   10, 10, 10, 10, 10, "2025-07-14 15:35:49", 10),
    @("luis",     7,  9,  5, 10, 10, "2025-07-14 15:37:47", 8.2),
    @("jessica", 10, 10, 10, 10, 10, "2025-07-14 15:51:55", 10)
)

$startRow = 24
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]

    $ws.Range("A$r").Value = $values[0]
    $ws.Range("B$r").Value = $values[1]
    $ws.Range("C$r").Value = $values[2]
    $ws.Range("D$r").Value = $values[3]
    $ws.Range("E$r").Value = $values[4]
    $ws.Range("F$r").Value = $values[5]
    $ws.Range("G$r").Value = $values[6]
    $ws.Range("H$r").Value = $values[7]
}
